# Update "Capital and operating costs" sheet:
#  - Recalculated capital cost figures (column C, rows 3-17)
#  - Reordered raw materials table (rows 21-27): Tridecane/CSL now come
#    before DAP/Glucose/Salt, "Process water" becomes its own single-row
#    "Raw materials" group after the "Wastewater" by-product row, and the
#    merged "A" label cells are adjusted accordingly.
#  - Updated operating-cost figures (columns C/D) that depend on the above.
#  - Updated maintenance/property-insurance figures (rows 35-36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First unmerge the existing merged label cells so every cell in the
#     raw-materials block can be written to independently. Borders are
#     temporarily cleared so that re-merging below doesn't make Excel
#     auto-split the thin border box around the new merged ranges (the
#     original formatting keeps one uniform bordered style for every row). ---
$ws.Range("A21:A23").UnMerge()
$ws.Range("A25:A28").UnMerge()
$ws.Range("A21:A28").Borders.LineStyle = 0

# --- Capital cost figures (column C, rows 3-17) ---
$ws.Range("C3").Value  = 61.38661436946352
$ws.Range("C4").Value  = 30.15271114172059
$ws.Range("C5").Value  = 2.455464574778541
$ws.Range("C6").Value  = 5.524795293251716
$ws.Range("C7").Value  = 2.762397646625858
$ws.Range("C8").Value  = 102.2819830258402
$ws.Range("C9").Value  = 10.22819830258402
$ws.Range("C10").Value = 10.22819830258402
$ws.Range("C11").Value = 20.45639660516805
$ws.Range("C12").Value = 30.68459490775206
$ws.Range("C13").Value = 10.22819830258402
$ws.Range("C14").Value = 81.82558642067215
$ws.Range("C15").Value = 184.1075694465124
$ws.Range("C16").Value = 9.205378472325618
$ws.Range("C17").Value = 193.312947918838

# --- Raw materials / by-products table (rows 21-28) ---
# Row 21: Tridecane (was DAP)
$ws.Range("B21").Value = "Tridecane"
$ws.Range("C21").Value = 878.1550799999999
$ws.Range("D21").Value = 0.002402640811183376

# Row 22: CSL (was Glucose)
$ws.Range("B22").Value = "CSL"
$ws.Range("C22").Value = 51.528108
$ws.Range("D22").Value = 0.5665466433993392

# Row 23: DAP (was Salt)
$ws.Range("B23").Value = "DAP"
$ws.Range("C23").Value = 895.3915949999999
$ws.Range("D23").Value = 1.181191106686225

# Row 24: Glucose (was Wastewater); no longer a group-start row
$ws.Range("A24").Value = $null
$ws.Range("B24").Value = "Glucose"
$ws.Range("C24").Value = 240.404025
$ws.Range("D24").Value = 144.7370872673688

# Row 25: Salt (was Process water); no longer a group-start row
$ws.Range("A25").Value = $null
$ws.Range("B25").Value = "Salt"
$ws.Range("C25").Value = 136.07775
$ws.Range("D25").Value = 4.230799213075651

# Row 26: Wastewater (was Tridecane); now starts "By-products and credits"
$ws.Range("A26").Value = "By-products and credits"
$ws.Range("B26").Value = "Wastewater"
$ws.Range("C26").Value = -1.962993797789321
$ws.Range("D26").Value = -6.266075981442587

# Row 27: Process water (was CSL); now starts its own "Raw materials" group
$ws.Range("A27").Value = "Raw materials"
$ws.Range("B27").Value = "Process water"
$ws.Range("C27").Value = 0.320236305
$ws.Range("D27").Value = 0.9954197635451122

# Row 28: Natural gas label/price unchanged, cost updated
$ws.Range("D28").Value = 2.868987199981105

# Row 29: Total variable operating cost
$ws.Range("D29").Value = 160.84850981631

# --- Re-merge the "A" label column with its new groupings, then restore
#     the thin-border box style that applies uniformly to the whole block. ---
$ws.Range("A21:A25").Merge()
$ws.Range("A27:A28").Merge()
$ws.Range("A21:A28").Borders.LineStyle = 1

# --- Maintenance / property insurance figures (rows 35-36) ---
$ws.Range("C35").Value = 1.841598431083905
$ws.Range("D35").Value = 1.767934493840549
$ws.Range("C36").Value = 0.4297063005862446
$ws.Range("D36").Value = 0.4125180485627948
